# add infix colon operator, start work on reference macros
# -> adds a new "Sheet5" worksheet with a small fruit price/count table and
#    a handful of INDEX()-based formulas that exercise multi-area refs and
#    the new ":" (range) operator applied to a function result.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet at the end of the workbook and name it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "Sheet5"

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Fruit"
$ws.Range("B1").Value = "Price"
$ws.Range("C1").Value = "Count"

# ---------------------------------------------------------------------
# 3. Data rows.
# ---------------------------------------------------------------------
$fruit = @("Apples","Bananas","Lemons","Oranges","Pears","Almonds","Cashews","Peanuts","Walnuts")
$price = @(0.69, 0.34, 0.55000000000000004, 0.25, 0.59, 2.8, 3.55, 1.25, 1.75)
$count = @(40, 38, 15, 25, 40, 10, 16, 20, 12)

for ($i = 0; $i -lt $fruit.Length; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $fruit[$i]
    $ws.Range("B$r").Value = $price[$i]
    $ws.Range("C$r").Value = $count[$i]
}

# ---------------------------------------------------------------------
# 4. Reference-macro formulas (multi-area refs, 0-col INDEX, ":" applied
#    to a function result).
# ---------------------------------------------------------------------
$ws.Range("A12").Formula = "=INDEX(A2:C6, 2, 3)"
$ws.Range("A13").Formula = "=INDEX((A1:C6, A8:C11), 2, 2, 2)"
$ws.Range("A14").Formula = "=SUM(INDEX(A1:C11, 0, 3, 1))"
$ws.Range("A15").Formula = "=SUM(B2:INDEX(A2:C6, 5, 2))"

# ---------------------------------------------------------------------
# 5. Formatting.
# ---------------------------------------------------------------------

# -- header row: bold Segoe UI, dark grey text on light grey fill,
#    medium light-grey border top+bottom, wrapped, left/center aligned.
$hdr = $ws.Range("A1:C1")
$hdr.Font.Name = "Segoe UI"
$hdr.Font.Bold = $true
$hdr.Font.Size = 19.2
$hdr.Font.Color = 3750201
$hdr.Interior.Color = 14342874
$hdr.Borders.Item(8).LineStyle = 1
$hdr.Borders.Item(8).Weight = -4138
$hdr.Borders.Item(8).Color = 13421772
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = -4138
$hdr.Borders.Item(9).Color = 13421772
$hdr.HorizontalAlignment = -4131
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# -- data rows (fruit name + count): Segoe UI, near-black text, light
#    fill, same medium borders, vertically centered + wrapped.
# (column A and column C are styled separately -- this runtime does not
#  propagate style writes across all areas of a multi-area/union range)
$labelsA = $ws.Range("A2:A10")
$labelsA.Font.Name = "Segoe UI"
$labelsA.Font.Size = 19.2
$labelsA.Font.Color = 1973790
$labelsA.Interior.Color = 16053492
$labelsA.Borders.Item(8).LineStyle = 1
$labelsA.Borders.Item(8).Weight = -4138
$labelsA.Borders.Item(8).Color = 13421772
$labelsA.Borders.Item(9).LineStyle = 1
$labelsA.Borders.Item(9).Weight = -4138
$labelsA.Borders.Item(9).Color = 13421772
$labelsA.VerticalAlignment = -4108
$labelsA.WrapText = $true

$labelsC = $ws.Range("C2:C10")
$labelsC.Font.Name = "Segoe UI"
$labelsC.Font.Size = 19.2
$labelsC.Font.Color = 1973790
$labelsC.Interior.Color = 16053492
$labelsC.Borders.Item(8).LineStyle = 1
$labelsC.Borders.Item(8).Weight = -4138
$labelsC.Borders.Item(8).Color = 13421772
$labelsC.Borders.Item(9).LineStyle = 1
$labelsC.Borders.Item(9).Weight = -4138
$labelsC.Borders.Item(9).Color = 13421772
$labelsC.VerticalAlignment = -4108
$labelsC.WrapText = $true

# -- price column: same look as the other data cells, plus currency format.
$prices = $ws.Range("B2:B10")
$prices.Font.Name = "Segoe UI"
$prices.Font.Size = 19.2
$prices.Font.Color = 1973790
$prices.Interior.Color = 16053492
$prices.Borders.Item(8).LineStyle = 1
$prices.Borders.Item(8).Weight = -4138
$prices.Borders.Item(8).Color = 13421772
$prices.Borders.Item(9).LineStyle = 1
$prices.Borders.Item(9).Weight = -4138
$prices.Borders.Item(9).Color = 13421772
$prices.VerticalAlignment = -4108
$prices.WrapText = $true
$prices.NumberFormat = """$""#,##0.00_);[Red](""$""#,##0.00)"

# -- A12 carries a plain wrap-text style; A13:A15 stay default.
$ws.Range("A12").WrapText = $true

# -- row heights / column widths (best effort; cosmetic only).
$ws.Rows.Item("1:10").RowHeight = 21
$ws.Rows.Item("12:15").RowHeight = 21
$ws.Columns.Item(1).ColumnWidth = 27.9
$ws.Columns.Item(2).ColumnWidth = 41.7
$ws.Columns.Item(3).ColumnWidth = 21.9

# -- portrait page orientation.
$ws.PageSetup.Orientation = 1

# -- leave the selection parked below the table, like the source file.
$ws.Range("A16").Select() | Out-Null

Write-Output "Sheet5 added."
